$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.764753333333333
$ws.Range("H2").Value = 14.29426
$ws.Range("I2").Value = 0.2966169987831952
$ws.Range("J2").Value = 0.2966169987831952
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3862596666666667
$ws.Range("N2").Value = 1.158779
$ws.Range("O2").Value = 0.05360826277999409
$ws.Range("P2").Value = 0.05360826277999409
$ws.Range("Q2").Value = 1.840432034282222
$ws.Range("R2").Value = 16.56388830854
$ws.Range("S2").Value = 0.01590112201578272
$ws.Range("T2").Value = 0.01590112201578272

$ws.Range("G3").Value = 4.764753333333333
$ws.Range("H3").Value = 14.29426
$ws.Range("I3").Value = 0.2966169987831952
$ws.Range("J3").Value = 0.2966169987831952
$ws.Range("O3").Value = 0.1630272174193556
$ws.Range("P3").Value = 0.1630272174193557
$ws.Range("Q3").Value = 5.596907973493334
$ws.Range("R3").Value = 50.37217176144
$ws.Range("S3").Value = 0.04835664395090471
$ws.Range("T3").Value = 0.04835664395090472

$ws.Range("G4").Value = 4.764753333333333
$ws.Range("H4").Value = 14.29426
$ws.Range("I4").Value = 0.2966169987831952
$ws.Range("J4").Value = 0.2966169987831952
$ws.Range("M4").Value = 5.644318666666667
$ws.Range("N4").Value = 16.932956
$ws.Range("O4").Value = 0.7833645198006502
$ws.Range("P4").Value = 0.7833645198006502
$ws.Range("Q4").Value = 26.89378618139555
$ws.Range("R4").Value = 242.04407563256
$ws.Range("S4").Value = 0.2323592328165078
$ws.Range("T4").Value = 0.2323592328165078

$ws.Range("I5").Value = 0.5337607564504776
$ws.Range("J5").Value = 0.5337607564504775
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3862596666666667
$ws.Range("N5").Value = 1.158779
$ws.Range("O5").Value = 0.05360826277999409
$ws.Range("P5").Value = 0.05360826277999409
$ws.Range("Q5").Value = 3.311847934690333
$ws.Range("R5").Value = 29.806631412213
$ws.Range("S5").Value = 0.02861398689344563
$ws.Range("T5").Value = 0.02861398689344562

$ws.Range("I6").Value = 0.5337607564504776
$ws.Range("J6").Value = 0.5337607564504775
$ws.Range("O6").Value = 0.1630272174193556
$ws.Range("P6").Value = 0.1630272174193557
$ws.Range("S6").Value = 0.08701753089177175
$ws.Range("T6").Value = 0.08701753089177174

$ws.Range("I7").Value = 0.5337607564504776
$ws.Range("J7").Value = 0.5337607564504775
$ws.Range("M7").Value = 5.644318666666667
$ws.Range("N7").Value = 16.932956
$ws.Range("O7").Value = 0.7833645198006502
$ws.Range("P7").Value = 0.7833645198006502
$ws.Range("Q7").Value = 48.39522925148134
$ws.Range("R7").Value = 435.557063263332
$ws.Range("S7").Value = 0.4181292386652602
$ws.Range("T7").Value = 0.4181292386652601

$ws.Range("G8").Value = 2.724753333333334
$ws.Range("H8").Value = 8.17426
$ws.Range("I8").Value = 0.1696222447663273
$ws.Range("J8").Value = 0.1696222447663273
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3862596666666667
$ws.Range("N8").Value = 1.158779
$ws.Range("O8").Value = 0.05360826277999409
$ws.Range("P8").Value = 0.05360826277999409
$ws.Range("Q8").Value = 1.052462314282222
$ws.Range("R8").Value = 9.47216082854
$ws.Range("S8").Value = 0.009093153870765751
$ws.Range("T8").Value = 0.009093153870765749

$ws.Range("G9").Value = 2.724753333333334
$ws.Range("H9").Value = 8.17426
$ws.Range("I9").Value = 0.1696222447663273
$ws.Range("J9").Value = 0.1696222447663273
$ws.Range("O9").Value = 0.1630272174193556
$ws.Range("P9").Value = 0.1630272174193557
$ws.Range("Q9").Value = 3.200626053493334
$ws.Range("R9").Value = 28.80563448144
$ws.Range("S9").Value = 0.0276530425766792
$ws.Range("T9").Value = 0.0276530425766792

$ws.Range("G10").Value = 2.724753333333334
$ws.Range("H10").Value = 8.17426
$ws.Range("I10").Value = 0.1696222447663273
$ws.Range("J10").Value = 0.1696222447663273
$ws.Range("M10").Value = 5.644318666666667
$ws.Range("N10").Value = 16.932956
$ws.Range("O10").Value = 0.7833645198006502
$ws.Range("P10").Value = 0.7833645198006502
$ws.Range("Q10").Value = 15.37937610139556
$ws.Range("R10").Value = 138.41438491256
$ws.Range("S10").Value = 0.1328760483188823
$ws.Range("T10").Value = 0.1328760483188823
